$d = $word.ActiveDocument

# 1. Narrative paragraph: "is amended to Possession of Marijuana Drug Paraphernalia"
#    -> "is amended to Driving Under Suspension FTA, Fines or Child Support"
$d.Content.Find.Execute(
    "is amended to Possession of Marijuana Drug Paraphernalia", $true, $false, $false, $false, $false,
    $true, 1, $false, "is amended to Driving Under Suspension FTA, Fines or Child Support", 2)

# 2. Table "Offense" row (caption line repeating the amendment).
#    Avoid touching the apostrophe in "Req'd" -- only replace the tail after "AMENDED to ".
$d.Content.Find.Execute(
    "AMENDED to Possession of Marijuana Drug Paraphernalia", $true, $false, $false, $false, $false,
    $true, 1, $false, "AMENDED to Driving Under Suspension FTA, Fines or Child Support", 2)

# 3. Table "Statute/Ord." row
$d.Content.Find.Execute(
    "2925.141(C) ", $true, $false, $false, $false, $false,
    $true, 1, $false, "4510.111 ", 2)

# 4. Table "Degree" row
$d.Content.Find.Execute(
    "Minor Misdemeanor", $true, $false, $false, $false, $false,
    $true, 1, $false, "Unclassified Misdemeanor", 2)

# 5. Table "Plea" row only (the first of the two "Guilty" cells -- "Finding" stays "Guilty")
$tbl = $d.Tables.Item(1)
foreach ($row in $tbl.Rows) {
    $labelCell = $row.Cells.Item(1)
    $labelText = $labelCell.Range.Text.TrimEnd([char]13, [char]7)
    if ($labelText -eq "Plea") {
        $valueCell = $row.Cells.Item(2)
        # Re-scope via Document.Range(start, end) -- Cell.Range.Find has been observed
        # to search/replace beyond the cell boundary in this runtime.
        $scoped = $d.Range($valueCell.Range.Start, $valueCell.Range.End)
        $scoped.Find.Execute(
            "Guilty", $true, $false, $false, $false, $false,
            $true, 0, $false, "No Contest", 2)
    }
}
